$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.040.75"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "3.079.95"
$ws.Range("E3").Value = "  -1.06%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("E6").Value = "  -2.89%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.077.37"
$ws.Range("E8").Value = "  -1.02%  "

$ws.Range("E9").Value = "  -1.26%  "

$ws.Range("E10").Value = "  -0.56%  "

$ws.Range("E11").Value = "  -1.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.472"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000242"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.75%  "

$ws.Range("E15").Value = "  -2.05%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "66.936.93"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "3.589.61"
$ws.Range("E17").Value = "  -1.07%  "

$ws.Range("E18").Value = "  -1.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.64%  "

$ws.Range("D20").Value = "3.081.99"
$ws.Range("E20").Value = "  -0.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "490.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.688"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.96%  "

$ws.Range("E27").Value = "  +2.97%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.67%  "

$ws.Range("E30").Value = "  -5.92%  "

$ws.Range("E31").Value = "  -1.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.92%  "

$ws.Range("E33").Value = "  -2.45%  "

$ws.Range("D34").Value = "0.0₃0908"
$ws.Range("E34").Value = "  -3.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.953"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "46.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.80%  "

$ws.Range("E39").Value = "  +0.83%  "

$ws.Range("E40").Value = "  -4.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.303"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.27%  "

$ws.Range("D43").Value = "2.774.60"
$ws.Range("E43").Value = "  -0.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "371.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0344"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "135.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.82%  "

$ws.Range("E49").Value = "  -1.93%  "

$ws.Range("E50").Value = "  -1.83%  "

$ws.Range("E51").Value = "  -1.63%  "
